$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 257 (shifts old 257..325 down to 259..327)
$ws.Range("A257:A258").EntireRow.Insert()

# New row 257: Primera, week of 44588
$ws.Cells.Item(257,1).Value = 4
$ws.Cells.Item(257,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(257,3).Value = "Los Lagos"
$ws.Cells.Item(257,4).Value = 44588
$ws.Cells.Item(257,5).Value = 10
$ws.Cells.Item(257,6).Value = "Fruta"
$ws.Cells.Item(257,7).Value = 100106
$ws.Cells.Item(257,8).Value = "Oleaginosos"
$ws.Cells.Item(257,9).Value = 100106002
$ws.Cells.Item(257,10).Value = "Palta"
$ws.Cells.Item(257,11).Value = "Hass"
$ws.Cells.Item(257,12).Value = "Primera"
$ws.Cells.Item(257,13).Value = 200
$ws.Cells.Item(257,14).Value = 4000
$ws.Cells.Item(257,15).Value = 4100
$ws.Cells.Item(257,16).Value = 4050
$ws.Cells.Item(257,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(257,18).Value = "Provincia de Quillota"
$ws.Cells.Item(257,19).Value = 4050
$ws.Cells.Item(257,20).Value = 1

# New row 258: Segunda, week of 44588
$ws.Cells.Item(258,1).Value = 4
$ws.Cells.Item(258,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(258,3).Value = "Los Lagos"
$ws.Cells.Item(258,4).Value = 44588
$ws.Cells.Item(258,5).Value = 10
$ws.Cells.Item(258,6).Value = "Fruta"
$ws.Cells.Item(258,7).Value = 100106
$ws.Cells.Item(258,8).Value = "Oleaginosos"
$ws.Cells.Item(258,9).Value = 100106002
$ws.Cells.Item(258,10).Value = "Palta"
$ws.Cells.Item(258,11).Value = "Hass"
$ws.Cells.Item(258,12).Value = "Segunda"
$ws.Cells.Item(258,13).Value = 100
$ws.Cells.Item(258,14).Value = 3500
$ws.Cells.Item(258,15).Value = 3500
$ws.Cells.Item(258,16).Value = 3500
$ws.Cells.Item(258,17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(258,18).Value = "Provincia de Quillota"
$ws.Cells.Item(258,19).Value = 3500
$ws.Cells.Item(258,20).Value = 1
